$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before existing row 106, shifting old rows 106-116 down to 109-119
# (this also grows the sheet dimension from R116 to R119 automatically).
$ws.Rows("106:108").Insert()

# Row 106 - new record
$ws.Range("A106").Value = 1
$ws.Range("B106").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C106").Value = "Arica y Parinacota"
$ws.Range("D106").Value = 45211
$ws.Range("E106").Value = 15
$ws.Range("F106").Value = 100112045
$ws.Range("G106").Value = "Zapallo"
$ws.Range("H106").Value = "Camote"
$ws.Range("I106").Value = "1a nueva(o)"
$ws.Range("J106").Value = 400
$ws.Range("K106").Value = 770
$ws.Range("L106").Value = 800
$ws.Range("M106").Value = 785
$ws.Range("N106").Value = "$/kilo (volumen en unidades)"
$ws.Range("O106").Value = "Perú"
$ws.Range("P106").Value = 785
$ws.Range("Q106").Value = 1
$ws.Range("R106").Value = "Hortaliza"

# Row 107 - new record
$ws.Range("A107").Value = 1
$ws.Range("B107").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C107").Value = "Arica y Parinacota"
$ws.Range("D107").Value = 45211
$ws.Range("E107").Value = 15
$ws.Range("F107").Value = 100112045
$ws.Range("G107").Value = "Zapallo"
$ws.Range("H107").Value = "Camote"
$ws.Range("I107").Value = "2a nueva(o)"
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 740
$ws.Range("L107").Value = 760
$ws.Range("M107").Value = 750
$ws.Range("N107").Value = "$/kilo (volumen en unidades)"
$ws.Range("O107").Value = "Perú"
$ws.Range("P107").Value = 750
$ws.Range("Q107").Value = 1
$ws.Range("R107").Value = "Hortaliza"

# Row 108 - new record
$ws.Range("A108").Value = 1
$ws.Range("B108").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C108").Value = "Arica y Parinacota"
$ws.Range("D108").Value = 45211
$ws.Range("E108").Value = 15
$ws.Range("F108").Value = 100112045
$ws.Range("G108").Value = "Zapallo"
$ws.Range("H108").Value = "Camote"
$ws.Range("I108").Value = "3a nueva (o)"
$ws.Range("J108").Value = 300
$ws.Range("K108").Value = 730
$ws.Range("L108").Value = 740
$ws.Range("M108").Value = 735
$ws.Range("N108").Value = "$/kilo (volumen en unidades)"
$ws.Range("O108").Value = "Perú"
$ws.Range("P108").Value = 735
$ws.Range("Q108").Value = 1
$ws.Range("R108").Value = "Hortaliza"
